# Apply the "Översikt BOXHOLM" update:
#  1) Bump the "Förändrad" (column C) date from 45188 to 45189 for every
#     existing data row (2..230).
#  2) Give row 230 an explicit 15pt custom row height (matches new rows).
#  3) Append two new records (rows 231 and 232).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Column C (Förändrad) 45188 -> 45189 for rows 2..230 -----------------
$ws.Range("C2:C230").Value = 45189

# --- 2) Row 230 gains an explicit custom row height -------------------------
$ws.Rows.Item(230).RowHeight = 15

# --- 3) New row 231: A 43807-2023 -------------------------------------------
$ws.Rows.Item(231).RowHeight = 15

$ws.Cells.Item(231, 1).Value = "A 43807-2023"

$ws.Cells.Item(231, 2).Value = 45187
$ws.Cells.Item(231, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(231, 3).Value = 45189
$ws.Cells.Item(231, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(231, 4).Value = "ÖSTERGÖTLANDS LÄN"
$ws.Cells.Item(231, 5).Value = "BOXHOLM"

$ws.Cells.Item(231, 7).Value = 2.2
$ws.Cells.Item(231, 8).Value = 0
$ws.Cells.Item(231, 9).Value = 0
$ws.Cells.Item(231, 10).Value = 0
$ws.Cells.Item(231, 11).Value = 0
$ws.Cells.Item(231, 12).Value = 0
$ws.Cells.Item(231, 13).Value = 0
$ws.Cells.Item(231, 14).Value = 0
$ws.Cells.Item(231, 15).Value = 0
$ws.Cells.Item(231, 16).Value = 0
$ws.Cells.Item(231, 17).Value = 0

$ws.Cells.Item(231, 18).Value = ""
$ws.Cells.Item(231, 18).WrapText = $true

# --- New row 232: A 43824-2023 ----------------------------------------------
$ws.Cells.Item(232, 1).Value = "A 43824-2023"

$ws.Cells.Item(232, 2).Value = 45187
$ws.Cells.Item(232, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(232, 3).Value = 45189
$ws.Cells.Item(232, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(232, 4).Value = "ÖSTERGÖTLANDS LÄN"
$ws.Cells.Item(232, 5).Value = "BOXHOLM"

$ws.Cells.Item(232, 7).Value = 1.6
$ws.Cells.Item(232, 8).Value = 0
$ws.Cells.Item(232, 9).Value = 0
$ws.Cells.Item(232, 10).Value = 0
$ws.Cells.Item(232, 11).Value = 0
$ws.Cells.Item(232, 12).Value = 0
$ws.Cells.Item(232, 13).Value = 0
$ws.Cells.Item(232, 14).Value = 0
$ws.Cells.Item(232, 15).Value = 0
$ws.Cells.Item(232, 16).Value = 0
$ws.Cells.Item(232, 17).Value = 0

$ws.Cells.Item(232, 18).Value = ""
$ws.Cells.Item(232, 18).WrapText = $true
